$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Countries list reorder: Venezuela moves up (now ranks above Taiwan & Reunion
#    by "Casos totales"), so rows 127-129 become Venezuela / Taiwan / Reunion
#    (row 130 "Sierra Leona" is unaffected).
$ws.Range("A127").Value = "Venezuela"
$ws.Range("A128").Value = "Taiwan"
$ws.Range("A129").Value = "Reunion"

# 2) Updated stats for the reordered rows.
$ws.Range("B127").Value = 440
$ws.Range("C127").Value = 0
$ws.Range("D127").Value = 220
$ws.Range("E127").Value = 210
$ws.Range("F127").Value = 2
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 10

$ws.Range("B128").Value = 440
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 375
$ws.Range("E128").Value = 58
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 7

$ws.Range("B129").Value = 439
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 354
$ws.Range("E129").Value = 85
$ws.Range("F129").Value = 4
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 0

# 3) Timestamp footer update.
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 08:05"

# 4) Germany ("Alemania") row refresh.
$ws.Range("D11").Value = 150300
$ws.Range("E11").Value = 15937
